$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 corresponds to 9c9bafea-51f4-4c8b-a552-1c000fb1ea13.md
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: row 3 is the same file - handback now in sync, clear stale error, refresh datetime
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-09-04 14:53:00"
$wsZhCn.Range("P3").Value = ""

# de-de sheet: row 3 is the same file - handback now in sync, clear stale error, refresh datetime
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-09-04 14:53:12"
$wsDeDe.Range("P3").Value = ""
